# Update cryptocurrency price/volume data per Thu May 9 2024 GitHub Actions run

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.567.43"
$ws.Range("E2").Value = "  -1.87%  "

$ws.Range("D3").Value = "2.997.15"
$ws.Range("E3").Value = "  -1.08%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.85"
$ws.Range("E5").Value = "  +2.09%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.99"
$ws.Range("E6").Value = "  -3.51%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.521"
$ws.Range("E8").Value = "  -0.83%  "

$ws.Range("D9").Value = "2.995.17"
$ws.Range("E9").Value = "  -1.18%  "

$ws.Range("E10").Value = "  -3.04%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.91"
$ws.Range("E11").Value = "  +4.42%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.462"
$ws.Range("E12").Value = "  +4.09%  "

$ws.Range("E13").Value = "  -1.45%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.27"
$ws.Range("E14").Value = "  -2.94%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.124"
$ws.Range("E15").Value = "  +2.53%  "

$ws.Range("D16").Value = "3.487.57"
$ws.Range("E16").Value = "  -1.26%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.01"
$ws.Range("E17").Value = "  -0.47%  "

$ws.Range("D18").Value = "61.513.31"
$ws.Range("E18").Value = "  -1.92%  "

$ws.Range("D19").Value = "2.993.92"
$ws.Range("E19").Value = "  -1.14%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "454.15"
$ws.Range("E20").Value = "  -3.05%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.04"
$ws.Range("E21").Value = "  -0.16%  "

$ws.Range("E22").Value = "  -0.50%  "

$ws.Range("E23").Value = "  -0.59%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.28"
$ws.Range("E24").Value = "  +1.47%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.19"
$ws.Range("E25").Value = "  -7.84%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.17"
$ws.Range("E26").Value = "  -2.10%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.45"
$ws.Range("E27").Value = "  +0.03%  "

$ws.Range("E28").Value = "  +0.09%  "

$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.67"
$ws.Range("E29").Value = "  +1.65%  "

$ws.Range("B30").Value = "FirstDigitalUSD"
$ws.Range("C30").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  -0.05%  "

$ws.Range("E31").Value = "  -2.72%  "

$ws.Range("E32").Value = "  -4.67%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.25"

$ws.Range("E34").Value = "  +0.04%  "

$ws.Range("D35").Value = "0.0₃0816"
$ws.Range("E35").Value = "  +2.15%  "

$ws.Range("E36").Value = "  -2.46%  "

$ws.Range("E37").Value = "  -0.47%  "

$ws.Range("E38").Value = "  -3.25%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.33"
$ws.Range("E40").Value = "  +0.16%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.122"
$ws.Range("E41").Value = "  +8.43%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.87"
$ws.Range("E42").Value = "  -3.44%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "400.00"
$ws.Range("E43").Value = "  -5.23%  "

$ws.Range("B44").Value = "Arweave"
$ws.Range("C44").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "39.28"
$ws.Range("E44").Value = "  +3.06%  "

$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0352"
$ws.Range("E45").Value = "  -0.86%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.267"
$ws.Range("E46").Value = "  -5.38%  "

$ws.Range("D47").Value = "2.716.07"
$ws.Range("E47").Value = "  -2.79%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "132.98"
$ws.Range("E48").Value = "  +2.43%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.107"
$ws.Range("E50").Value = "  -0.56%  "

$ws.Range("E51").Value = "  +1.02%  "

